$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6933.3335
$ws.Range("I40").Value = 6000
$ws.Range("K40").Value = 6000
$ws.Range("M40").Value = -5825

$ws.Range("H69").Value = 20000
$ws.Range("I69").Value = 20000
$ws.Range("K69").Value = 60000
$ws.Range("M69").Value = -59126

$ws.Range("H72").Value = 20000
$ws.Range("I72").Value = 20000
$ws.Range("K72").Value = 180000
$ws.Range("M72").Value = -175632

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H137").Value = 2036.4333
$ws.Range("I137").Value = 1576.238
$ws.Range("K137").Value = 4728.714
$ws.Range("M137").Value = -2178.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 13785.857
$ws.Range("J13").Value = 13785.857
$ws.Range("L13").Value = 13785.857
$ws.Range("N13").Value = -14073.857

$ws.Range("H14").Value = 16996.666
$ws.Range("I14").Value = 16996.666
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 16996.666
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -16821.666
$ws.Range("N14").ClearContents()

$ws.Range("H16").Value = 50503
$ws.Range("I16").Value = 50503
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 50503
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -50216
$ws.Range("N16").ClearContents()

$ws.Range("H88").Value = 2178.9285
$ws.Range("I88").Value = 295.6
$ws.Range("J88").Value = 3225.2222
$ws.Range("K88").Value = 295.6
$ws.Range("L88").Value = 3225.2222
$ws.Range("M88").Value = 110.4
$ws.Range("N88").Value = -4037.2222

$ws.Range("H91").Value = 2178.9285
$ws.Range("I91").Value = 295.6
$ws.Range("J91").Value = 3225.2222
$ws.Range("K91").Value = 295.6
$ws.Range("L91").Value = 3225.2222
$ws.Range("M91").Value = 1108.4
$ws.Range("N91").Value = -6033.2222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1561.125
$ws.Range("I20").Value = 1555.5714
$ws.Range("K20").Value = 1555.5714
$ws.Range("M20").Value = -1308.5714

$ws.Range("H99").Value = 3699.7646
$ws.Range("I99").Value = 4071.2144
$ws.Range("J99").Value = 1966.3334
$ws.Range("K99").Value = 4071.2144
$ws.Range("L99").Value = 1966.3334
$ws.Range("M99").Value = -2573.2144
$ws.Range("N99").Value = -4962.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 82.5
$ws.Range("I7").Value = 76.666664
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 76.666664
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = 36.333336
$ws.Range("N7").Value = -326

$ws.Range("H31").Value = 1674
$ws.Range("I31").Value = 1674
$ws.Range("K31").Value = 1674
$ws.Range("M31").Value = -1379

$ws.Range("H34").Value = 1674
$ws.Range("I34").Value = 1674
$ws.Range("K34").Value = 1674
$ws.Range("M34").Value = -1472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 94
$ws.Range("J12").Value = 125
$ws.Range("L12").Value = 375
$ws.Range("N12").Value = -721

$ws.Range("H33").Value = 227.14285
$ws.Range("J33").Value = 596.5
$ws.Range("L33").Value = 3579
$ws.Range("N33").Value = -4145

$ws.Range("H55").Value = 2060
$ws.Range("J55").Value = 2426.25
$ws.Range("L55").Value = 7278.75
$ws.Range("N55").Value = -7632.75

$ws.Range("H59").Value = 500
$ws.Range("I59").Value = 500
$ws.Range("K59").Value = 1500
$ws.Range("M59").Value = -960

$ws.Range("H75").Value = 1905
$ws.Range("I75").Value = 4800
$ws.Range("J75").Value = 457.5
$ws.Range("K75").Value = 14400
$ws.Range("L75").Value = 1372.5
$ws.Range("M75").Value = -13402
$ws.Range("N75").Value = -3368.5

$ws.Range("H78").Value = 1905
$ws.Range("I78").Value = 4800
$ws.Range("J78").Value = 457.5
$ws.Range("K78").Value = 43200
$ws.Range("L78").Value = 4117.5
$ws.Range("M78").Value = -38208
$ws.Range("N78").Value = -14101.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 272.91306
$ws.Range("I2").Value = 341.25
$ws.Range("J2").Value = 116.71429
$ws.Range("K2").Value = 341.25
$ws.Range("L2").Value = 116.71429
$ws.Range("M2").Value = -228.25
$ws.Range("N2").Value = -342.71429

$ws.Range("H11").Value = 335000.8
$ws.Range("J11").Value = 425001.34
$ws.Range("L11").Value = 425001.34
$ws.Range("N11").Value = -425279.34

$ws.Range("H70").Value = 30305910
$ws.Range("I70").Value = 37040000
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 37040000
$ws.Range("L70").Value = 2500
$ws.Range("M70").Value = -37039730
$ws.Range("N70").Value = -3040

$ws.Range("H73").Value = 30305910
$ws.Range("I73").Value = 37040000
$ws.Range("J73").Value = 2500
$ws.Range("K73").Value = 37040000
$ws.Range("L73").Value = 2500
$ws.Range("M73").Value = -37039064
$ws.Range("N73").Value = -4372

$ws.Range("H123").Value = 25333.666
$ws.Range("J123").Value = 25333.666
$ws.Range("L123").Value = 25333.666
$ws.Range("N123").Value = -30233.666

$ws.Range("H126").Value = 10239.8
$ws.Range("I126").Value = 10239.8
$ws.Range("K126").Value = 30719.4
$ws.Range("M126").Value = -28249.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3533.5454
$ws.Range("I46").Value = 2976.3635
$ws.Range("K46").Value = 2976.3635
$ws.Range("M46").Value = -2788.3635

$ws.Range("H61").Value = 2614.1
$ws.Range("I61").Value = 2614.1
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2614.1
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2412.1
$ws.Range("N61").ClearContents()

$ws.Range("H62").Value = 28613
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248

$ws.Range("H65").Value = 28613
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

$ws.Range("H93").Value = 3583.3333
$ws.Range("I93").Value = 3583.3333
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 3583.3333
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -2335.3333
$ws.Range("N93").ClearContents()

$ws.Range("H108").Value = 37500
$ws.Range("J108").Value = 37500
$ws.Range("L108").Value = 37500
$ws.Range("N108").Value = -45180

$ws.Range("H113").Value = 2614.1
$ws.Range("I113").Value = 2614.1
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2614.1
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -444.0999999999999
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 17174
$ws.Range("I23").Value = 17174
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 17174
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -16945
$ws.Range("N23").ClearContents()

$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41872

$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129360

$ws.Range("H81").Value = 1242.6666
$ws.Range("J81").Value = 1116
$ws.Range("L81").Value = 2232
$ws.Range("N81").Value = -4354

$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30766

$ws.Range("H84").Value = 1242.6666
$ws.Range("J84").Value = 1116
$ws.Range("L84").Value = 11160
$ws.Range("N84").Value = -21768

$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32652

$ws.Range("H107").Value = 2486.4211
$ws.Range("I107").Value = 2932.818
$ws.Range("J107").Value = 1872.625
$ws.Range("K107").Value = 8798.454000000002
$ws.Range("L107").Value = 5617.875
$ws.Range("M107").Value = -6878.454000000002
$ws.Range("N107").Value = -9457.875
